# Auto-applies the cryptos.xlsx price/volume update described in the commit
# "Updated cryptos list on Wed Mar 15 10:48:04 UTC 2023 with GitHub Actions".
#
# Column D (Price) and E (Volume(1h)) are stored as literal text (inlineStr)
# in the workbook, e.g. "24.689.52" or "  +0.90%  ". Plain numeric-looking
# strings (like "1.005") would otherwise be auto-coerced to a number by Excel's
# COM layer, so those are entered with a leading apostrophe to force text,
# exactly as a human typing the figures into a pre-formatted text column would.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple rows: only Price (D) and Volume(1h) (E) change ---
$ws.Cells.Item(2, 4).Value = "24.769.00"
$ws.Cells.Item(2, 5).Value = "  +1.21%  "
$ws.Cells.Item(3, 4).Value = "1.696.82"
$ws.Cells.Item(3, 5).Value = "  +1.04%  "
$ws.Cells.Item(4, 4).Value = "'1.005"
$ws.Cells.Item(4, 5).Value = "  +0.22%  "
$ws.Cells.Item(5, 4).Value = "'311.02"
$ws.Cells.Item(5, 5).Value = "  +1.64%  "
$ws.Cells.Item(6, 4).Value = "'1.001"
$ws.Cells.Item(6, 5).Value = "  +0.41%  "
$ws.Cells.Item(7, 4).Value = "'0.3718"
$ws.Cells.Item(7, 5).Value = "  +0.86%  "
$ws.Cells.Item(8, 4).Value = "'48.82"
$ws.Cells.Item(8, 5).Value = "  +1.85%  "
$ws.Cells.Item(9, 4).Value = "'0.3409"
$ws.Cells.Item(9, 5).Value = "  -0.47%  "
$ws.Cells.Item(10, 4).Value = "'1.205"
$ws.Cells.Item(10, 5).Value = "  +3.82%  "
$ws.Cells.Item(11, 4).Value = "'0.07428"
$ws.Cells.Item(11, 5).Value = "  +2.88%  "
$ws.Cells.Item(12, 4).Value = "'1.004"
$ws.Cells.Item(12, 5).Value = "  +0.52%  "
$ws.Cells.Item(13, 4).Value = "'6.293"
$ws.Cells.Item(13, 5).Value = "  +2.95%  "
$ws.Cells.Item(14, 4).Value = "'20.85"
$ws.Cells.Item(14, 5).Value = "  +3.50%  "
$ws.Cells.Item(15, 4).Value = "'6.943"
$ws.Cells.Item(15, 5).Value = "  +3.34%  "
$ws.Cells.Item(16, 4).Value = "1.690.30"
$ws.Cells.Item(16, 5).Value = "  +0.82%  "
$ws.Cells.Item(17, 4).Value = "'0.00001118"
$ws.Cells.Item(17, 5).Value = "  +1.55%  "
$ws.Cells.Item(18, 4).Value = "'0.06701"
$ws.Cells.Item(18, 5).Value = "  +0.76%  "
$ws.Cells.Item(19, 4).Value = "'1.001"
$ws.Cells.Item(19, 5).Value = "  +0.41%  "
$ws.Cells.Item(20, 4).Value = "'82.97"
$ws.Cells.Item(20, 5).Value = "  +2.79%  "
$ws.Cells.Item(21, 4).Value = "'17.10"
$ws.Cells.Item(21, 5).Value = "  +4.01%  "
$ws.Cells.Item(22, 4).Value = "'6.297"
$ws.Cells.Item(22, 5).Value = "  +3.46%  "
$ws.Cells.Item(23, 4).Value = "'12.87"
$ws.Cells.Item(23, 5).Value = "  +6.36%  "
$ws.Cells.Item(24, 4).Value = "24.729.89"
$ws.Cells.Item(24, 5).Value = "  +1.35%  "
$ws.Cells.Item(25, 4).Value = "'2.447"
$ws.Cells.Item(25, 5).Value = "  +1.22%  "
$ws.Cells.Item(26, 4).Value = "'2.755"
$ws.Cells.Item(26, 5).Value = "  +3.91%  "
$ws.Cells.Item(27, 4).Value = "'20.15"
$ws.Cells.Item(27, 5).Value = "  +3.76%  "
$ws.Cells.Item(28, 4).Value = "'148.25"
$ws.Cells.Item(28, 5).Value = "  -3.17%  "
$ws.Cells.Item(29, 4).Value = "'131.37"
$ws.Cells.Item(29, 5).Value = "  +3.37%  "
$ws.Cells.Item(30, 4).Value = "1.878.94"
$ws.Cells.Item(30, 5).Value = "  +0.88%  "
$ws.Cells.Item(31, 4).Value = "'1.239"
$ws.Cells.Item(31, 5).Value = "  +26.97%  "
$ws.Cells.Item(32, 4).Value = "'6.694"
$ws.Cells.Item(32, 5).Value = "  +7.01%  "
$ws.Cells.Item(33, 4).Value = "'4.225"
$ws.Cells.Item(33, 5).Value = "  +4.88%  "
$ws.Cells.Item(34, 4).Value = "'13.56"
$ws.Cells.Item(34, 5).Value = "  +9.77%  "
$ws.Cells.Item(37, 4).Value = "'5.514"
$ws.Cells.Item(37, 5).Value = "  +3.79%  "
$ws.Cells.Item(38, 4).Value = "'0.06580"
$ws.Cells.Item(38, 5).Value = "  +3.23%  "
$ws.Cells.Item(39, 4).Value = "'9.024"
$ws.Cells.Item(39, 5).Value = "  +4.19%  "
$ws.Cells.Item(40, 4).Value = "'0.02390"
$ws.Cells.Item(40, 5).Value = "  +3.37%  "
$ws.Cells.Item(41, 4).Value = "'0.2201"
$ws.Cells.Item(41, 5).Value = "  +5.40%  "
$ws.Cells.Item(43, 4).Value = "'0.6370"
$ws.Cells.Item(43, 5).Value = "  +4.53%  "
$ws.Cells.Item(45, 4).Value = "'13.64"
$ws.Cells.Item(45, 5).Value = "  +5.71%  "
$ws.Cells.Item(48, 4).Value = "'2.100"
$ws.Cells.Item(48, 5).Value = "  +4.50%  "
$ws.Cells.Item(49, 4).Value = "'128.08"
$ws.Cells.Item(49, 5).Value = "  +1.99%  "
$ws.Cells.Item(50, 4).Value = "'0.07229"
$ws.Cells.Item(50, 5).Value = "  +0.93%  "
$ws.Cells.Item(51, 4).Value = "'78.95"
$ws.Cells.Item(51, 5).Value = "  +4.31%  "

# --- Rows whose Coin/Link (B/C) were reordered, with new D/E values ---
$ws.Cells.Item(35, 2).Value = "WEMIXTOKEN"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(35, 4).Value = "'1.752"
$ws.Cells.Item(35, 5).Value = "  +3.36%  "
$ws.Cells.Item(36, 2).Value = "Stellar"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(36, 4).Value = "'0.08662"
$ws.Cells.Item(36, 5).Value = "  +2.78%  "
$ws.Cells.Item(46, 2).Value = "PancakeSwap"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(46, 4).Value = "'3.811"
$ws.Cells.Item(46, 5).Value = "  +1.30%  "
$ws.Cells.Item(47, 2).Value = "Decentraland"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Cells.Item(47, 4).Value = "'0.6038"
$ws.Cells.Item(47, 5).Value = "  +2.74%  "

# --- Rows where only Volume(1h) (E) changes ---
$ws.Cells.Item(42, 5).Value = "  +0.86%  "
$ws.Cells.Item(44, 5).Value = "  +0.49%  "
